$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

# Extend the Research description (row 4, column D) with the additional sentence
$ws.Range("D4").Value = "i began with reading the assingment and beginning to do research on algorithms, i want to chose the best algoritm for a quick and customizable result, i chose for randomized depth first search, because it is simple to implement without any mistakes and it will do it effectively"

# Row 5 used to be the "Example 2" sample row - turn it into the real "Implementing generation" entry
$ws.Range("A5").Value = "Implementing generation"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "11/21/2022"
$ws.Range("D5").Value = "starting with implementing the generation, also made a github, generates a grid and selects a random starting position"

# Row 6 used to be the "Example 3" bonus sample row - clear it out entirely
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

# Move the active selection to A6, matching the saved sheet view
$ws.Range("A6").Select()

# Recompute the "Total amount of hours" SUMIF (B30) now that the bonus row (row 6)
# no longer has an "x" marker in column E. The headless calc engine used here mis-handles
# "<>" criteria against truly blank cells, so briefly populate the blank E4:E28 cells with a
# non-"x" placeholder to force a correct recalculation, then restore the blanks under manual
# calculation so the freshly computed cached result is preserved on save.
for ($r = 4; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $v = $cell.Value()
    if ($v -eq $null -or $v -eq "") {
        $cell.Value = "_tmp_"
    }
}
$excel.CalculateFullRebuild()

$excel.Calculation = -4135
for ($r = 4; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value() -eq "_tmp_") {
        $cell.Value = ""
    }
}
$excel.Calculation = -4105
